# feat: add 2022-Q3 data
#
# - Insert a new "2022-Q3" sheet right after "总计", pushing the other
#   quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) one tab to the right.
# - Add a new row to "总计" for the 2022-Q3 totals, shifting the existing
#   quarter rows down by one.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)      # "总计"
$q2     = $wb.Worksheets.Item(2)      # "2022-Q2" (existing, used only as a style template)

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q3" and renumber the rest
# ---------------------------------------------------------------------
$totals.Rows.Item(2).Insert()

# Copy the formatting of the (now shifted down) first data row onto the
# newly inserted blank row so the new row matches the existing look.
$totals.Range("A3:D3").Copy()
$totals.Range("A2:D2").PasteSpecial(-4122)   # xlPasteFormats

$totalsData = @(
  @(0, "2022-Q3", 2, 0.01),
  @(1, "2022-Q2", 2, 0.01),
  @(2, "2022-Q1", 2, 0.01),
  @(3, "2021-Q4", 3, 0.01),
  @(4, "2021-Q3", 2, 0.01)
)
for ($i = 0; $i -lt $totalsData.Count; $i++) {
  $r = $i + 2
  $row = $totalsData[$i]
  $totals.Cells.Item($r, 1).Value = $row[0]
  $totals.Cells.Item($r, 2).Value = $row[1]
  $totals.Cells.Item($r, 3).Value = $row[2]
  $totals.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, placed right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totals)
$q3.Name = "2022-Q3"

# Header row styling to match the other quarter sheets
$q3Header = $q3.Range("B1:H1")
$q3Header.Font.Bold = $true
$q3Header.Borders.LineStyle = 1
$q3Header.HorizontalAlignment = -4108
$q3Header.VerticalAlignment = -4160

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
  $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3IdCol = $q3.Range("A2:A3")
$q3IdCol.Font.Bold = $true
$q3IdCol.Borders.LineStyle = 1
$q3IdCol.HorizontalAlignment = -4108
$q3IdCol.VerticalAlignment = -4160

# Force columns B..G to stay text (fund codes with leading zeros, and the
# percentage-looking figures are stored as text throughout this workbook).
$q3.Range("B2:G3").NumberFormat = "@"

$q3Rows = @(
  @(0, "010343", "华宝英国富时100指数（QDII）A", "0.13", "92.85", "2.81", "0.0037", 9),
  @(1, "010344", "华宝英国富时100指数（QDII）C", "0.08", "92.85", "2.81", "0.0022", 9)
)
for ($i = 0; $i -lt $q3Rows.Count; $i++) {
  $r = $i + 2
  $row = $q3Rows[$i]
  $q3.Cells.Item($r, 1).Value = $row[0]
  $q3.Cells.Item($r, 2).Value = $row[1]
  $q3.Cells.Item($r, 3).Value = $row[2]
  $q3.Cells.Item($r, 4).Value = $row[3]
  $q3.Cells.Item($r, 5).Value = $row[4]
  $q3.Cells.Item($r, 6).Value = $row[5]
  $q3.Cells.Item($r, 7).Value = $row[6]
  $q3.Cells.Item($r, 8).Value = $row[7]
}
